# Update cryptocurrency price and 1h-volume-change figures.
#
# The Price/Volume columns are stored as literal text (e.g. "38.152.77",
# "  +2.97%  ") rather than numbers, so each cell is written as Text to avoid
# Excel's automatic numeric conversion (e.g. "230.26" silently becoming the
# number 230.26 instead of staying the string "230.26"), and ClearFormats is
# used right after so the cell keeps its original (default) style instead of
# picking up the temporary "@" text format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $range = $ws.Range($address)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue "D2" "38.157.27"
Set-TextValue "E2" "  +3.03%  "
Set-TextValue "D3" "2.060.37"
Set-TextValue "E3" "  +2.66%  "
Set-TextValue "D5" "230.26"
Set-TextValue "E5" "  +2.20%  "
Set-TextValue "D6" "0.618"
Set-TextValue "E6" "  +1.87%  "
Set-TextValue "D7" "59.48"
Set-TextValue "E7" "  +8.34%  "
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "D9" "0.386"
Set-TextValue "E9" "  +3.65%  "
Set-TextValue "D10" "0.0810"
Set-TextValue "E10" "  +4.50%  "
Set-TextValue "E11" "  +2.89%  "
Set-TextValue "E12" "  +5.98%  "
Set-TextValue "D13" "2.364.53"
Set-TextValue "E13" "  +2.50%  "
Set-TextValue "D14" "21.30"
Set-TextValue "E14" "  +8.40%  "
Set-TextValue "E15" "  +3.17%  "
Set-TextValue "E16" "  +2.58%  "
Set-TextValue "D17" "2.053.55"
Set-TextValue "E17" "  +1.09%  "
Set-TextValue "D18" "38.063.33"
Set-TextValue "E18" "  +3.00%  "
Set-TextValue "E19" "  +1.66%  "
Set-TextValue "D20" "69.89"
Set-TextValue "E20" "  +2.45%  "
Set-TextValue "E21" "  +3.29%  "
Set-TextValue "D22" "225.25"
Set-TextValue "E22" "  +0.85%  "
Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  +0.01%  "
Set-TextValue "D24" "2.44"
Set-TextValue "E24" "  +0.46%  "
Set-TextValue "E25" "  +4.86%  "
Set-TextValue "D26" "9.30"
Set-TextValue "E26" "  +4.89%  "
Set-TextValue "D27" "166.28"
Set-TextValue "E27" "  +1.33%  "
Set-TextValue "D28" "0.134"
Set-TextValue "E28" "  +8.65%  "
Set-TextValue "D29" "19.06"
Set-TextValue "E29" "  +2.78%  "
Set-TextValue "E30" "  +3.35%  "
Set-TextValue "E31" "  +2.80%  "
Set-TextValue "E32" "  +3.61%  "
Set-TextValue "E33" "  +2.97%  "
Set-TextValue "E34" "  +10.50%  "
Set-TextValue "E35" "  +1.87%  "
Set-TextValue "E36" "  +1.75%  "
Set-TextValue "D37" "6.17"
Set-TextValue "E37" "  +16.15%  "
Set-TextValue "E38" "  +5.90%  "
Set-TextValue "E39" "  +0.14%  "
Set-TextValue "D40" "1.534.83"
Set-TextValue "E40" "  +5.72%  "
Set-TextValue "D41" "98.35"
Set-TextValue "E41" "  +4.23%  "
Set-TextValue "E42" "  +3.04%  "
Set-TextValue "D43" "16.93"
Set-TextValue "E43" "  +6.87%  "
Set-TextValue "E44" "  +4.40%  "
Set-TextValue "E45" "  +2.03%  "
Set-TextValue "E46" "  +1.59%  "
Set-TextValue "D47" "4.17"
Set-TextValue "E47" "  +4.43%  "
Set-TextValue "E48" "  +3.29%  "
Set-TextValue "D49" "2.98"
Set-TextValue "E49" "  +3.31%  "
Set-TextValue "D50" "7.12"
Set-TextValue "E50" "  +0.74%  "
Set-TextValue "D51" "2.251.85"
Set-TextValue "E51" "  +2.61%  "
